$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK_RegisteredDatas")

$ws.Range("A3").Value = "VO22SSU"
$ws.Range("B3").Value = "Pasquale"
$ws.Range("C3").Value = "Johnston"
$ws.Range("D3").Value = "Carola"
$ws.Range("E3").Value = "Wyman"
$ws.Range("F3").Value = "kareem.denesik@hotmail.com"
$ws.Range("G3").Value = "test1234"
$ws.Range("H3").Value = "'7043203860"
$ws.Range("H3").Style = "Normal"
